# Atualização - RDD 13 e Copas
# Adds the "Rodada 13" column (N) with scores for each team.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new round, copying the style (bold, bordered, centered)
# used by the other "Rodada n" header cells.
$ws.Range("N1").Value = "Rodada 13"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# New round scores per team (rows 2-21)
$ws.Range("N2").Value = 70.2099609375
$ws.Range("N3").Value = 67.91015625
$ws.Range("N4").Value = 102.2099609375
$ws.Range("N5").Value = 108.509765625
$ws.Range("N6").Value = 100.75
$ws.Range("N7").Value = 118.740234375
$ws.Range("N8").Value = 107.64990234375
$ws.Range("N9").Value = 121.2099609375
$ws.Range("N10").Value = 132.2099609375
$ws.Range("N11").Value = 118.740234375
$ws.Range("N12").Value = 117.10986328125
$ws.Range("N13").Value = 105.91015625
$ws.Range("N14").Value = 133.4404296875
$ws.Range("N15").Value = 113.669921875
$ws.Range("N16").Value = 121.0400390625
$ws.Range("N17").Value = 117.509765625
$ws.Range("N18").Value = 111.91015625
$ws.Range("N19").Value = 79.41015625
$ws.Range("N20").Value = 130.2099609375
$ws.Range("N21").Value = 133.0498046875
